# Update "想去人数" (interest count) values in column F for several convention
# entries that appear on both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row -> new value
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value  = 295
$wsExpo.Range("F6").Value  = 147
$wsExpo.Range("F9").Value  = 2074
$wsExpo.Range("F10").Value = 363
$wsExpo.Range("F11").Value = 5000
$wsExpo.Range("F12").Value = 102
$wsExpo.Range("F13").Value = 348

# Sheet "全部类型" (all types) - same events, different row offsets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 295
$wsAll.Range("F8").Value  = 147
$wsAll.Range("F13").Value = 2074
$wsAll.Range("F14").Value = 363
$wsAll.Range("F15").Value = 5000
$wsAll.Range("F16").Value = 102
$wsAll.Range("F17").Value = 348

$wb.Save()
